# Applies the "Jurisdiction" sheet restructuring described by the commit:
#   "Uploading single feature layers for all (1) zoned and (2) unzoned
#    jurisdictions analyzed to date."
#
# The Village-specific (unzoned-town-only split) data is removed from the
# Jurisdiction sheet, a dedicated "County" row is introduced, and the
# jurisdiction name is normalized from "Town of Ludlow" to "Ludlow Town".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Jurisdiction")

# Insert a new row 2 for the "County" field (pushes everything else down by one).
$ws.Rows.Item(2).Insert()

$ws.Range("A2").Value = "County"
$ws.Range("B2").Value = "Windsor"

# Jurisdiction Name: "Town of Ludlow" -> "Ludlow Town"
$ws.Range("B1").Value = "Ludlow Town"

# "Pages in Zoning Code" row (was row 6, now row 7): keep only the Town page
# number (as a plain number) and drop the Village-specific note.
$ws.Range("B7").Value = 98
$ws.Range("C7").ClearContents()

# "Link to Zoning Code" row (was row 7, now row 8): drop the Village zoning
# link entirely (value + hyperlink), keep the Town link in B8 untouched.
foreach ($hl in @($ws.Hyperlinks)) {
    if ($hl.Range.Address() -eq "`$C`$8") {
        $hl.Delete()
    }
}
$ws.Range("C8").ClearContents()

# "Effective Date Zoning Text" row (was row 12, now row 13): drop the Village
# effective date, keep the Town date in B13 untouched.
$ws.Range("C13").ClearContents()

$wb.Save()
